$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cost and quantity for the HB100 line item; the Total Cost
# formula (E3 = C3*D3) recalculates automatically.
$ws.Range("C3").Value = 228.9
$ws.Range("D3").Value = 2

# Move the active selection to D12 (as last left by the editing session).
$ws.Range("D12").Select()
